# "simplifying dFBA flux units to facilitate composability"
#
# On the "dFBA objectives" sheet, split the old single "Coefficient units"
# column into two columns: a new "Reaction rate units" column (holding the
# units of the net reaction rate, e.g. "s^-1") inserted immediately before
# the existing "Coefficient units" column (which keeps its old value, e.g. "s").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dFBA objectives")

# Insert a new blank column at F, pushing "Coefficient units" (and everything
# to its right) one column to the right.
$ws.Columns.Item(6).Insert() | Out-Null

# Populate the new column's header and the single data row.
$ws.Cells.Item(1, 6).Value = "Reaction rate units"
$ws.Cells.Item(2, 6).Value = "s^-1"

# Incidental: the reviewer also selected the "Flux min" column on the
# Reactions sheet while looking at this change.
$wsReactions = $wb.Worksheets.Item("Reactions")
$wsReactions.Range("F:F").Select() | Out-Null

# Restore focus/selection to the dFBA objectives sheet.
$ws.Range("F1:F2").Select() | Out-Null
$ws.Activate() | Out-Null
